# Edit script: duplicate rows 34-41 twice (to rows 42-49 and 50-57)
# and update the uuid (column G) for all data rows (2-57) to the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUuid = "9d7fe935-9ba4-4289-a505-a9d49b45ad9f"
$newUuid = "ddb71f0e-ca75-4b40-9ae2-33afa81c43ba"

# Source block is rows 34-41 (8 rows), columns A-H
$srcFirstRow = 34
$srcLastRow = 41
$numRows = $srcLastRow - $srcFirstRow + 1  # 8

# Capture the source block values/formulas once before we overwrite the uuid,
# since we'll be writing the new uuid into the same cells.
$srcRange = $ws.Range("A$srcFirstRow`:H$srcLastRow")

# Duplicate the block twice: first copy goes to rows 42-49, second to rows 50-57
for ($copyIndex = 1; $copyIndex -le 2; $copyIndex++) {
    $destFirstRow = $srcLastRow + ($copyIndex - 1) * $numRows + 1
    $destRange = $ws.Range("A$destFirstRow`:H" + ($destFirstRow + $numRows - 1))
    $srcRange.Copy($destRange)
}

# Now update column G (uuid) for every data row from 2 through 57 (the new last row)
$lastRow = 57
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldUuid) {
        $cell.Value = $newUuid
    }
}
